$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D19").Value = "2016-03-03 11:04:35"
$wsZhCn.Range("D20").Value = "2016-03-03 11:04:35"
$wsZhCn.Range("G19").Value = "2016-03-03 11:05:25"
$wsZhCn.Range("G20").Value = "2016-03-03 11:05:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D19").Value = "2016-03-03 11:04:47"
$wsDeDe.Range("D20").Value = "2016-03-03 11:04:47"
$wsDeDe.Range("G19").Value = "2016-03-03 11:05:47"
$wsDeDe.Range("G20").Value = "2016-03-03 11:05:47"
